$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows to append (id_registro, id_usuario, accion, fecha_accion)
$rows = @(
    @(45, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-18 15:39:29"),
    @(46, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-18 15:40:47"),
    @(47, "nico@gmail.com", "Actualizó pedido #1 de 'pagado' a 'pendiente'", "2026-02-18 15:41:20"),
    @(48, "nico@gmail.com", "Actualizó pedido #1 de 'pendiente' a 'enviado'", "2026-02-18 15:41:23"),
    @(49, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-18 17:34:40")
)

$startRow = 46
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $data[3]
}
